# Generate Report for Handoff
# Regenerated the handoff XLIFF files for the 4 "Ready for handoff" rows
# (rows 4-7) on both the zh-cn and de-de sheets: their Priority moved from
# "low" to "ht" and their Latest Handoff Datetime was refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-21 06:39:27"
$overview.Range("G5").Value = "2016-08-21 06:39:27"
$overview.Range("G6").Value = "2016-08-21 06:39:27"
$overview.Range("G7").Value = "2016-08-21 06:39:27"

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("E4:E7").Value = "ht"
$zh.Range("H4").Value = "2016-08-21 06:39:23"
$zh.Range("H5").Value = "2016-08-21 06:39:23"
$zh.Range("H6").Value = "2016-08-21 06:39:23"
$zh.Range("H7").Value = "2016-08-21 06:39:23"

$de = $wb.Worksheets.Item("de-de")
$de.Range("E4:E7").Value = "ht"
$de.Range("H4").Value = "2016-08-21 06:39:27"
$de.Range("H5").Value = "2016-08-21 06:39:27"
$de.Range("H6").Value = "2016-08-21 06:39:27"
$de.Range("H7").Value = "2016-08-21 06:39:27"
